$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.708.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.513.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.512.58"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.123"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.381"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.104.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.506.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.684.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.573"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.652.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -3.92%  "
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.98%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.513.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0811"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.817"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.45%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.376.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0266"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.59%  "
